$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage (matches original inlineStr cells)
# by setting NumberFormat to Text before assignment, then resetting the
# style afterwards so no stray style index is left on the cell.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Update Price (D) and Volume(1h) (E) columns per the data refresh ---
Set-TextValue $ws.Range("D2") "65.855.57"
Set-TextValue $ws.Range("E2") "  -1.90%  "
Set-TextValue $ws.Range("D3") "3.436.15"
Set-TextValue $ws.Range("E3") "  -0.47%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "582.88"
Set-TextValue $ws.Range("E5") "  -0.36%  "
Set-TextValue $ws.Range("D6") "173.34"
Set-TextValue $ws.Range("E6") "  -1.63%  "
Set-TextValue $ws.Range("E7") "  -0.05%  "
Set-TextValue $ws.Range("D8") "0.605"
Set-TextValue $ws.Range("E8") "  -0.60%  "
Set-TextValue $ws.Range("D9") "3.432.02"
Set-TextValue $ws.Range("E9") "  -0.56%  "
Set-TextValue $ws.Range("E10") "  -3.27%  "
Set-TextValue $ws.Range("D11") "6.92"
Set-TextValue $ws.Range("E11") "  -0.22%  "
Set-TextValue $ws.Range("D12") "0.408"
Set-TextValue $ws.Range("E12") "  -3.64%  "
Set-TextValue $ws.Range("D13") "4.033.96"
Set-TextValue $ws.Range("E13") "  -0.75%  "
Set-TextValue $ws.Range("E14") "  +1.33%  "
Set-TextValue $ws.Range("D15") "28.67"
Set-TextValue $ws.Range("E15") "  -9.21%  "
Set-TextValue $ws.Range("D16") "65.977.11"
Set-TextValue $ws.Range("E16") "  -1.75%  "
Set-TextValue $ws.Range("D17") "0.0000170"
Set-TextValue $ws.Range("E17") "  -2.02%  "
Set-TextValue $ws.Range("D18") "3.449.02"
Set-TextValue $ws.Range("E18") "  +0.20%  "
Set-TextValue $ws.Range("D19") "5.92"
Set-TextValue $ws.Range("E19") "  -2.37%  "
Set-TextValue $ws.Range("D20") "13.81"
Set-TextValue $ws.Range("E20") "  -0.42%  "
Set-TextValue $ws.Range("D21") "368.73"
Set-TextValue $ws.Range("E21") "  -2.73%  "
Set-TextValue $ws.Range("D22") "7.64"
Set-TextValue $ws.Range("E22") "  -2.14%  "
Set-TextValue $ws.Range("D23") "72.37"
Set-TextValue $ws.Range("E23") "  +1.06%  "
Set-TextValue $ws.Range("E24") "  +0.37%  "
Set-TextValue $ws.Range("D25") "0.533"
Set-TextValue $ws.Range("E25") "  +0.44%  "
Set-TextValue $ws.Range("D26") "0.0000121"
Set-TextValue $ws.Range("E26") "  +0.88%  "
Set-TextValue $ws.Range("D27") "9.72"
Set-TextValue $ws.Range("E27") "  -1.67%  "
Set-TextValue $ws.Range("D28") "0.176"
Set-TextValue $ws.Range("E28") "  +1.60%  "
Set-TextValue $ws.Range("E29") "  -0.02%  "
Set-TextValue $ws.Range("D30") "23.62"
Set-TextValue $ws.Range("E30") "  -1.38%  "
Set-TextValue $ws.Range("D31") "5.72"
Set-TextValue $ws.Range("E31") "  -3.93%  "
Set-TextValue $ws.Range("D32") "1.97"
Set-TextValue $ws.Range("E32") "  -2.37%  "
Set-TextValue $ws.Range("E33") "  +0.02%  "
Set-TextValue $ws.Range("D34") "1.28"
Set-TextValue $ws.Range("E34") "  -5.91%  "
Set-TextValue $ws.Range("D35") "6.99"
Set-TextValue $ws.Range("E35") "  -2.78%  "
Set-TextValue $ws.Range("E36") "  -0.70%  "
Set-TextValue $ws.Range("D37") "160.50"
Set-TextValue $ws.Range("E37") "  +0.24%  "
Set-TextValue $ws.Range("D40") "1.76"
Set-TextValue $ws.Range("E40") "  -2.67%  "
Set-TextValue $ws.Range("D41") "2.59"
Set-TextValue $ws.Range("E41") "  -1.56%  "
Set-TextValue $ws.Range("D42") "2.768.51"
Set-TextValue $ws.Range("E42") "  +2.68%  "
Set-TextValue $ws.Range("E43") "  -2.58%  "
Set-TextValue $ws.Range("D44") "4.45"
Set-TextValue $ws.Range("E44") "  -0.52%  "
Set-TextValue $ws.Range("D45") "0.0680"
Set-TextValue $ws.Range("E45") "  -2.26%  "
Set-TextValue $ws.Range("D46") "40.12"
Set-TextValue $ws.Range("E46") "  -2.30%  "
Set-TextValue $ws.Range("E47") "  -3.94%  "
Set-TextValue $ws.Range("D48") "0.0289"
Set-TextValue $ws.Range("E48") "  -1.51%  "
Set-TextValue $ws.Range("D49") "323.80"
Set-TextValue $ws.Range("E49") "  +0.72%  "
Set-TextValue $ws.Range("E50") "  -1.55%  "
Set-TextValue $ws.Range("D51") "6.24"
Set-TextValue $ws.Range("E51") "  +0.18%  "

# --- Rows 38/39: EnergySwap and Mantle swap positions, with refreshed data ---
Set-TextValue $ws.Range("B38") "Mantle"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D38") "0.877"
Set-TextValue $ws.Range("E38") "  -0.01%  "

Set-TextValue $ws.Range("B39") "EnergySwap"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D39") "28.65"
Set-TextValue $ws.Range("E39") "  +6.57%  "
